$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H6").Value = 1632.1
$ws_ALC.Range("I6").Value = 590.125
$ws_ALC.Range("K6").Value = 1770.375
$ws_ALC.Range("M6").Value = -1658.375
$ws_ALC.Range("H9").Value = 516.2222
$ws_ALC.Range("I9").Value = 289.5
$ws_ALC.Range("J9").Value = 1513.8
$ws_ALC.Range("K9").Value = 289.5
$ws_ALC.Range("L9").Value = 1513.8
$ws_ALC.Range("M9").Value = -120.5
$ws_ALC.Range("N9").Value = -1851.8
$ws_ALC.Range("H17").Value = 3475.862
$ws_ALC.Range("I17").Value = 0
$ws_ALC.Range("J17").Value = 3475.862
$ws_ALC.Range("K17").Value = 0
$ws_ALC.Range("L17").Value = 10427.586
$ws_ALC.Range("M17").ClearContents()
$ws_ALC.Range("N17").Value = -10763.586
$ws_ALC.Range("H43").Value = 2772.1072
$ws_ALC.Range("I43").Value = 2531.3
$ws_ALC.Range("J43").Value = 3374.125
$ws_ALC.Range("K43").Value = 2531.3
$ws_ALC.Range("L43").Value = 3374.125
$ws_ALC.Range("M43").Value = -2462.3
$ws_ALC.Range("N43").Value = -3512.125
$ws_ALC.Range("H107").Value = 777.25
$ws_ALC.Range("I107").Value = 796.73334
$ws_ALC.Range("J107").Value = 485
$ws_ALC.Range("K107").Value = 796.73334
$ws_ALC.Range("L107").Value = 485
$ws_ALC.Range("M107").Value = 1123.26666
$ws_ALC.Range("N107").Value = -4325
$ws_ALC.Range("H121").Value = 2785.0625
$ws_ALC.Range("J121").Value = 2785.0625
$ws_ALC.Range("L121").Value = 8355.1875
$ws_ALC.Range("N121").Value = -11849.1875
$ws_ALC.Range("H127").Value = 1977.2727
$ws_ALC.Range("I127").Value = 797.1177
$ws_ALC.Range("K127").Value = 2391.3531
$ws_ALC.Range("M127").Value = 2568.6469
$ws_ALC.Range("H132").Value = 4257.091
$ws_ALC.Range("I132").Value = 4221.7144
$ws_ALC.Range("J132").Value = 5000
$ws_ALC.Range("K132").Value = 12665.1432
$ws_ALC.Range("L132").Value = 15000
$ws_ALC.Range("M132").Value = -10135.1432
$ws_ALC.Range("N132").Value = -20060
$ws_ALC.Range("H137").Value = 3280.535
$ws_ALC.Range("I137").Value = 2901.8857
$ws_ALC.Range("J137").Value = 4937.125
$ws_ALC.Range("K137").Value = 8705.6571
$ws_ALC.Range("L137").Value = 14811.375
$ws_ALC.Range("M137").Value = -6155.6571
$ws_ALC.Range("N137").Value = -19911.375
$ws_ALC.Range("H138").Value = 4730.1465
$ws_ALC.Range("I138").Value = 5523.6665
$ws_ALC.Range("J138").Value = 4531.7666
$ws_ALC.Range("K138").Value = 16570.9995
$ws_ALC.Range("L138").Value = 13595.2998
$ws_ALC.Range("M138").Value = -11430.9995
$ws_ALC.Range("N138").Value = -23875.2998
$ws_ALC.Range("H141").Value = 2658.75
$ws_ALC.Range("I141").Value = 1487.8667
$ws_ALC.Range("J141").Value = 4610.222
$ws_ALC.Range("K141").Value = 4463.6001
$ws_ALC.Range("L141").Value = 13830.666
$ws_ALC.Range("M141").Value = 716.3999000000003
$ws_ALC.Range("N141").Value = -24190.666
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H2").Value = 2620.303
$ws_ARM.Range("I2").Value = 2774
$ws_ARM.Range("J2").Value = 2140
$ws_ARM.Range("K2").Value = 2774
$ws_ARM.Range("L2").Value = 2140
$ws_ARM.Range("M2").Value = -2661
$ws_ARM.Range("N2").Value = -2366
$ws_ARM.Range("H32").Value = 16920.385
$ws_ARM.Range("I32").Value = 12514.275
$ws_ARM.Range("K32").Value = 12514.275
$ws_ARM.Range("M32").Value = -12227.275
$ws_ARM.Range("H45").Value = 3790.6
$ws_ARM.Range("I45").Value = 3367.75
$ws_ARM.Range("K45").Value = 3367.75
$ws_ARM.Range("M45").Value = -2990.75
$ws_ARM.Range("H49").Value = 18590.773
$ws_ARM.Range("I49").Value = 17777.555
$ws_ARM.Range("J49").Value = 19882.354
$ws_ARM.Range("K49").Value = 17777.555
$ws_ARM.Range("L49").Value = 19882.354
$ws_ARM.Range("M49").Value = -17517.555
$ws_ARM.Range("N49").Value = -20402.354
$ws_ARM.Range("H53").Value = 39999
$ws_ARM.Range("J53").Value = 39999
$ws_ARM.Range("L53").Value = 39999
$ws_ARM.Range("N53").Value = -41363
$ws_ARM.Range("H59").Value = 25999.5
$ws_ARM.Range("I59").Value = 0
$ws_ARM.Range("J59").Value = 25999.5
$ws_ARM.Range("K59").Value = 0
$ws_ARM.Range("L59").Value = 25999.5
$ws_ARM.Range("M59").ClearContents()
$ws_ARM.Range("N59").Value = -27607.5
$ws_ARM.Range("H60").Value = 16000
$ws_ARM.Range("I60").Value = 0
$ws_ARM.Range("J60").Value = 16000
$ws_ARM.Range("K60").Value = 0
$ws_ARM.Range("L60").Value = 16000
$ws_ARM.Range("M60").ClearContents()
$ws_ARM.Range("N60").Value = -17466
$ws_ARM.Range("H61").Value = 8444.950999999999
$ws_ARM.Range("I61").Value = 7979.162
$ws_ARM.Range("J61").Value = 12753.5
$ws_ARM.Range("K61").Value = 7979.162
$ws_ARM.Range("L61").Value = 12753.5
$ws_ARM.Range("M61").Value = -7767.162
$ws_ARM.Range("N61").Value = -13177.5
$ws_ARM.Range("H92").Value = 20000
$ws_ARM.Range("J92").Value = 20000
$ws_ARM.Range("L92").Value = 20000
$ws_ARM.Range("N92").Value = -24992
$ws_ARM.Range("H97").Value = 1115.6666
$ws_ARM.Range("I97").Value = 0
$ws_ARM.Range("K97").Value = 0
$ws_ARM.Range("M97").ClearContents()
$ws_ARM.Range("H110").Value = 1506
$ws_ARM.Range("I110").Value = 1454.2222
$ws_ARM.Range("J110").Value = 1661.3334
$ws_ARM.Range("K110").Value = 1454.2222
$ws_ARM.Range("L110").Value = 1661.3334
$ws_ARM.Range("M110").Value = 590.7778000000001
$ws_ARM.Range("N110").Value = -5751.3334
$ws_ARM.Range("H116").Value = 2620.303
$ws_ARM.Range("I116").Value = 2774
$ws_ARM.Range("J116").Value = 2140
$ws_ARM.Range("K116").Value = 2774
$ws_ARM.Range("L116").Value = 2140
$ws_ARM.Range("M116").Value = -480
$ws_ARM.Range("N116").Value = -6728
$ws_ARM.Range("H122").Value = 3246
$ws_ARM.Range("I122").Value = 2609.4614
$ws_ARM.Range("J122").Value = 6004.3335
$ws_ARM.Range("K122").Value = 7828.3842
$ws_ARM.Range("L122").Value = 18013.0005
$ws_ARM.Range("M122").Value = -5378.3842
$ws_ARM.Range("N122").Value = -22913.0005
$ws_ARM.Range("H129").Value = 54780
$ws_ARM.Range("J129").Value = 54780
$ws_ARM.Range("L129").Value = 54780
$ws_ARM.Range("N129").Value = -64780
$ws_ARM.Range("H132").Value = 10985.714
$ws_ARM.Range("I132").Value = 4633.3335
$ws_ARM.Range("K132").Value = 13900.0005
$ws_ARM.Range("M132").Value = -11370.0005
$ws_ARM.Range("H136").Value = 8444.950999999999
$ws_ARM.Range("I136").Value = 7979.162
$ws_ARM.Range("J136").Value = 12753.5
$ws_ARM.Range("K136").Value = 23937.486
$ws_ARM.Range("L136").Value = 38260.5
$ws_ARM.Range("M136").Value = -21387.486
$ws_ARM.Range("N136").Value = -43360.5
$ws_ARM.Range("H139").Value = 128714.25
$ws_ARM.Range("J139").Value = 128714.25
$ws_ARM.Range("L139").Value = 128714.25
$ws_ARM.Range("N139").Value = -138994.25
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H3").Value = 2620.303
$ws_BSM.Range("I3").Value = 2774
$ws_BSM.Range("J3").Value = 2140
$ws_BSM.Range("K3").Value = 2774
$ws_BSM.Range("L3").Value = 2140
$ws_BSM.Range("M3").Value = -2660
$ws_BSM.Range("N3").Value = -2368
$ws_BSM.Range("H105").Value = 1588.9423
$ws_BSM.Range("J105").Value = 586.1177
$ws_BSM.Range("L105").Value = 586.1177
$ws_BSM.Range("N105").Value = -4080.1177
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H4").Value = 26974.75
$ws_CRP.Range("I4").Value = 21949.5
$ws_CRP.Range("J4").Value = 32000
$ws_CRP.Range("K4").Value = 21949.5
$ws_CRP.Range("L4").Value = 32000
$ws_CRP.Range("M4").Value = -21837.5
$ws_CRP.Range("N4").Value = -32224
$ws_CRP.Range("H7").Value = 2627.8372
$ws_CRP.Range("I7").Value = 2538.4644
$ws_CRP.Range("K7").Value = 2538.4644
$ws_CRP.Range("M7").Value = -2425.4644
$ws_CRP.Range("H16").Value = 1941
$ws_CRP.Range("I16").Value = 954.7778
$ws_CRP.Range("K16").Value = 954.7778
$ws_CRP.Range("M16").Value = -667.7778
$ws_CRP.Range("H22").Value = 2283
$ws_CRP.Range("I22").Value = 1224.75
$ws_CRP.Range("J22").Value = 2812.125
$ws_CRP.Range("K22").Value = 1224.75
$ws_CRP.Range("L22").Value = 2812.125
$ws_CRP.Range("M22").Value = -874.75
$ws_CRP.Range("N22").Value = -3512.125
$ws_CRP.Range("H31").Value = 6213.982
$ws_CRP.Range("I31").Value = 3258.2144
$ws_CRP.Range("J31").Value = 15763.385
$ws_CRP.Range("K31").Value = 3258.2144
$ws_CRP.Range("L31").Value = 15763.385
$ws_CRP.Range("M31").Value = -2963.2144
$ws_CRP.Range("N31").Value = -16353.385
$ws_CRP.Range("H34").Value = 6213.982
$ws_CRP.Range("I34").Value = 3258.2144
$ws_CRP.Range("J34").Value = 15763.385
$ws_CRP.Range("K34").Value = 3258.2144
$ws_CRP.Range("L34").Value = 15763.385
$ws_CRP.Range("M34").Value = -3056.2144
$ws_CRP.Range("N34").Value = -16167.385
$ws_CRP.Range("H99").Value = 28573.25
$ws_CRP.Range("I99").Value = 34764.668
$ws_CRP.Range("J99").Value = 9999
$ws_CRP.Range("K99").Value = 34764.668
$ws_CRP.Range("L99").Value = 9999
$ws_CRP.Range("M99").Value = -33266.668
$ws_CRP.Range("N99").Value = -12995
$ws_CRP.Range("H107").Value = 2258.5
$ws_CRP.Range("I107").Value = 2208
$ws_CRP.Range("J107").Value = 2372.125
$ws_CRP.Range("K107").Value = 2208
$ws_CRP.Range("L107").Value = 2372.125
$ws_CRP.Range("M107").Value = -288
$ws_CRP.Range("N107").Value = -6212.125
$ws_CRP.Range("H113").Value = 1941
$ws_CRP.Range("I113").Value = 954.7778
$ws_CRP.Range("K113").Value = 954.7778
$ws_CRP.Range("M113").Value = 1215.2222
$ws_CRP.Range("H122").Value = 4408.9414
$ws_CRP.Range("J122").Value = 4746.5386
$ws_CRP.Range("L122").Value = 14239.6158
$ws_CRP.Range("N122").Value = -19139.6158
$ws_CRP.Range("H126").Value = 28573.25
$ws_CRP.Range("I126").Value = 34764.668
$ws_CRP.Range("J126").Value = 9999
$ws_CRP.Range("K126").Value = 104294.004
$ws_CRP.Range("L126").Value = 29997
$ws_CRP.Range("M126").Value = -101824.004
$ws_CRP.Range("N126").Value = -34937
$ws_CRP.Range("H132").Value = 13042.4
$ws_CRP.Range("I132").Value = 11803
$ws_CRP.Range("J132").Value = 18000
$ws_CRP.Range("K132").Value = 35409
$ws_CRP.Range("L132").Value = 54000
$ws_CRP.Range("M132").Value = -32879
$ws_CRP.Range("N132").Value = -59060
$ws_CRP.Range("H134").Value = 3698.8064
$ws_CRP.Range("I134").Value = 3081.6667
$ws_CRP.Range("J134").Value = 5814.7144
$ws_CRP.Range("K134").Value = 9245.000100000001
$ws_CRP.Range("L134").Value = 17444.1432
$ws_CRP.Range("M134").Value = -6710.000100000001
$ws_CRP.Range("N134").Value = -22514.1432
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H44").Value = 1430.75
$ws_CUL.Range("I44").Value = 111.5
$ws_CUL.Range("K44").Value = 334.5
$ws_CUL.Range("M44").Value = 63.5
$ws_CUL.Range("H58").Value = 11116.2
$ws_CUL.Range("I58").Value = 11116.2
$ws_CUL.Range("K58").Value = 33348.60000000001
$ws_CUL.Range("M58").Value = -33220.60000000001
$ws_CUL.Range("H70").Value = 10036.25
$ws_CUL.Range("I70").Value = 4058
$ws_CUL.Range("K70").Value = 12174
$ws_CUL.Range("M70").Value = -11859
$ws_CUL.Range("H73").Value = 10036.25
$ws_CUL.Range("I73").Value = 4058
$ws_CUL.Range("K73").Value = 12174
$ws_CUL.Range("M73").Value = -11082
$ws_CUL.Range("H80").Value = 3000
$ws_CUL.Range("I80").Value = 1000
$ws_CUL.Range("J80").Value = 5000
$ws_CUL.Range("K80").Value = 3000
$ws_CUL.Range("L80").Value = 15000
$ws_CUL.Range("M80").Value = -2064
$ws_CUL.Range("N80").Value = -16872
$ws_CUL.Range("H83").Value = 3000
$ws_CUL.Range("I83").Value = 1000
$ws_CUL.Range("J83").Value = 5000
$ws_CUL.Range("K83").Value = 9000
$ws_CUL.Range("L83").Value = 45000
$ws_CUL.Range("M83").Value = -4320
$ws_CUL.Range("N83").Value = -54360
$ws_CUL.Range("H104").Value = 3396
$ws_CUL.Range("J104").Value = 3396
$ws_CUL.Range("L104").Value = 10188
$ws_CUL.Range("N104").Value = -15430
$ws_CUL.Range("H107").Value = 496.33334
$ws_CUL.Range("I107").Value = 307.16666
$ws_CUL.Range("J107").Value = 590.9167
$ws_CUL.Range("K107").Value = 921.4999799999999
$ws_CUL.Range("L107").Value = 1772.7501
$ws_CUL.Range("M107").Value = 998.5000200000001
$ws_CUL.Range("N107").Value = -5612.7501
$ws_CUL.Range("H120").Value = 15037.9375
$ws_CUL.Range("I120").Value = 9075.700000000001
$ws_CUL.Range("K120").Value = 27227.1
$ws_CUL.Range("M120").Value = -22389.1
$ws_CUL.Range("H122").Value = 3649
$ws_CUL.Range("I122").Value = 3649
$ws_CUL.Range("K122").Value = 32841
$ws_CUL.Range("M122").Value = -30391
$ws_CUL.Range("H130").Value = 2052.3333
$ws_CUL.Range("I130").Value = 2052.3333
$ws_CUL.Range("K130").Value = 6156.999899999999
$ws_CUL.Range("M130").Value = -1136.999899999999
$ws_CUL.Range("H138").Value = 4111.857
$ws_CUL.Range("I138").Value = 3997.6
$ws_CUL.Range("J138").Value = 4397.5
$ws_CUL.Range("K138").Value = 11992.8
$ws_CUL.Range("L138").Value = 13192.5
$ws_CUL.Range("M138").Value = -6852.799999999999
$ws_CUL.Range("N138").Value = -23472.5
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H2").Value = 3621.842
$ws_GSM.Range("J2").Value = 2167
$ws_GSM.Range("L2").Value = 2167
$ws_GSM.Range("N2").Value = -2393
$ws_GSM.Range("H26").Value = 23332.334
$ws_GSM.Range("J26").Value = 29999
$ws_GSM.Range("L26").Value = 29999
$ws_GSM.Range("N26").Value = -30559
$ws_GSM.Range("H50").Value = 23332.334
$ws_GSM.Range("J50").Value = 29999
$ws_GSM.Range("L50").Value = 29999
$ws_GSM.Range("N50").Value = -30995
$ws_GSM.Range("H113").Value = 2548.7693
$ws_GSM.Range("I113").Value = 2126.2222
$ws_GSM.Range("K113").Value = 2126.2222
$ws_GSM.Range("M113").Value = 43.77779999999984
$ws_GSM.Range("H122").Value = 6642.7856
$ws_GSM.Range("I122").Value = 4200.15
$ws_GSM.Range("J122").Value = 12749.375
$ws_GSM.Range("K122").Value = 12600.45
$ws_GSM.Range("L122").Value = 38248.125
$ws_GSM.Range("M122").Value = -10150.45
$ws_GSM.Range("N122").Value = -43148.125
$ws_GSM.Range("H126").Value = 9612.286
$ws_GSM.Range("I126").Value = 11547.777
$ws_GSM.Range("J126").Value = 8160.6665
$ws_GSM.Range("K126").Value = 34643.331
$ws_GSM.Range("L126").Value = 24481.9995
$ws_GSM.Range("M126").Value = -32173.331
$ws_GSM.Range("N126").Value = -29421.9995
$ws_GSM.Range("H132").Value = 5840.8
$ws_GSM.Range("I132").Value = 5811.6055
$ws_GSM.Range("J132").Value = 5999.2856
$ws_GSM.Range("K132").Value = 17434.8165
$ws_GSM.Range("L132").Value = 17997.8568
$ws_GSM.Range("M132").Value = -14904.8165
$ws_GSM.Range("N132").Value = -23057.8568
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H2").Value = 40000
$ws_LTW.Range("J2").Value = 40000
$ws_LTW.Range("L2").Value = 40000
$ws_LTW.Range("N2").Value = -40224
$ws_LTW.Range("H7").Value = 502500
$ws_LTW.Range("I7").Value = 1000000
$ws_LTW.Range("K7").Value = 1000000
$ws_LTW.Range("M7").Value = -999888
$ws_LTW.Range("H18").Value = 46250
$ws_LTW.Range("H20").Value = 9420
$ws_LTW.Range("I20").Value = 7304
$ws_LTW.Range("K20").Value = 7304
$ws_LTW.Range("M20").Value = -7078
$ws_LTW.Range("H22").Value = 4326.6665
$ws_LTW.Range("I22").Value = 4326.6665
$ws_LTW.Range("K22").Value = 4326.6665
$ws_LTW.Range("M22").Value = -4031.6665
$ws_LTW.Range("H27").Value = 4326.6665
$ws_LTW.Range("I27").Value = 4326.6665
$ws_LTW.Range("K27").Value = 4326.6665
$ws_LTW.Range("M27").Value = -4219.6665
$ws_LTW.Range("H40").Value = 4289.25
$ws_LTW.Range("I40").Value = 3534.5
$ws_LTW.Range("J40").Value = 6553.5
$ws_LTW.Range("K40").Value = 3534.5
$ws_LTW.Range("L40").Value = 6553.5
$ws_LTW.Range("M40").Value = -3398.5
$ws_LTW.Range("N40").Value = -6825.5
$ws_LTW.Range("H46").Value = 3390.9473
$ws_LTW.Range("I46").Value = 3868.375
$ws_LTW.Range("K46").Value = 3868.375
$ws_LTW.Range("M46").Value = -3680.375
$ws_LTW.Range("H82").Value = 4765.1
$ws_LTW.Range("I82").Value = 4758.5
$ws_LTW.Range("J82").Value = 4775
$ws_LTW.Range("K82").Value = 4758.5
$ws_LTW.Range("L82").Value = 4775
$ws_LTW.Range("M82").Value = -4397.5
$ws_LTW.Range("N82").Value = -5497
$ws_LTW.Range("H85").Value = 4765.1
$ws_LTW.Range("I85").Value = 4758.5
$ws_LTW.Range("J85").Value = 4775
$ws_LTW.Range("K85").Value = 4758.5
$ws_LTW.Range("L85").Value = 4775
$ws_LTW.Range("M85").Value = -3510.5
$ws_LTW.Range("N85").Value = -7271
$ws_LTW.Range("H93").Value = 5633.3335
$ws_LTW.Range("I93").Value = 3450
$ws_LTW.Range("K93").Value = 3450
$ws_LTW.Range("M93").Value = -2202
$ws_LTW.Range("H126").Value = 502500
$ws_LTW.Range("I126").Value = 1000000
$ws_LTW.Range("K126").Value = 3000000
$ws_LTW.Range("M126").Value = -2997530
$ws_LTW.Range("H132").Value = 11325
$ws_LTW.Range("I132").Value = 0
$ws_LTW.Range("K132").Value = 0
$ws_LTW.Range("M132").ClearContents()
$ws_LTW.Range("H136").Value = 7212.17
$ws_LTW.Range("I136").Value = 4167.5884
$ws_LTW.Range("J136").Value = 8649.888999999999
$ws_LTW.Range("K136").Value = 12502.7652
$ws_LTW.Range("L136").Value = 25949.667
$ws_LTW.Range("M136").Value = -9952.765199999998
$ws_LTW.Range("N136").Value = -31049.667
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H2").Value = 16500
$ws_WVR.Range("I2").Value = 16500
$ws_WVR.Range("K2").Value = 16500
$ws_WVR.Range("M2").Value = -16388
$ws_WVR.Range("H6").Value = 12502.5
$ws_WVR.Range("I6").Value = 1005
$ws_WVR.Range("J6").Value = 24000
$ws_WVR.Range("K6").Value = 1005
$ws_WVR.Range("L6").Value = 24000
$ws_WVR.Range("M6").Value = -890
$ws_WVR.Range("N6").Value = -24230
$ws_WVR.Range("H15").Value = 18706.076
$ws_WVR.Range("I15").Value = 15691.286
$ws_WVR.Range("J15").Value = 22223.334
$ws_WVR.Range("K15").Value = 15691.286
$ws_WVR.Range("L15").Value = 22223.334
$ws_WVR.Range("M15").Value = -15403.286
$ws_WVR.Range("N15").Value = -22799.334
$ws_WVR.Range("H54").Value = 37855.43
$ws_WVR.Range("J54").Value = 32998
$ws_WVR.Range("L54").Value = 32998
$ws_WVR.Range("N54").Value = -34038
$ws_WVR.Range("H64").Value = 125752
$ws_WVR.Range("J64").Value = 125752
$ws_WVR.Range("L64").Value = 125752
$ws_WVR.Range("N64").Value = -126248
$ws_WVR.Range("H67").Value = 125752
$ws_WVR.Range("J67").Value = 125752
$ws_WVR.Range("L67").Value = 125752
$ws_WVR.Range("N67").Value = -127468
$ws_WVR.Range("H96").Value = 23500
$ws_WVR.Range("I96").Value = 0
$ws_WVR.Range("J96").Value = 23500
$ws_WVR.Range("K96").Value = 0
$ws_WVR.Range("L96").Value = 23500
$ws_WVR.Range("M96").ClearContents()
$ws_WVR.Range("N96").Value = -26246
$ws_WVR.Range("H112").Value = 22623.5
$ws_WVR.Range("J112").Value = 22623.5
$ws_WVR.Range("L112").Value = 22623.5
$ws_WVR.Range("N112").Value = -25577.5
$ws_WVR.Range("H126").Value = 4762.3076
$ws_WVR.Range("I126").Value = 3892.1
$ws_WVR.Range("K126").Value = 11676.3
$ws_WVR.Range("M126").Value = -9206.299999999999
$ws_WVR.Range("H132").Value = 5377.4062
$ws_WVR.Range("I132").Value = 2917.1365
$ws_WVR.Range("J132").Value = 10790
$ws_WVR.Range("K132").Value = 8751.4095
$ws_WVR.Range("L132").Value = 32370
$ws_WVR.Range("M132").Value = -6221.4095
$ws_WVR.Range("N132").Value = -37430
$ws_WVR.Range("H135").Value = 0
$ws_WVR.Range("J135").Value = 0
$ws_WVR.Range("L135").Value = 0
$ws_WVR.Range("N135").ClearContents()
